$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells to reflect renamed codes / names
$ws.Range("B2").Value = "1040-20%A"
$ws.Range("I2").Value = "external.A"
$ws.Range("J2").Value = "ritenuta"

$ws.Range("B3").Value = "1040-23%R"
$ws.Range("I3").Value = "external.R"
$ws.Range("J3").Value = "ritenuta"

$ws.Range("C4").Value = "Enasarco 17% su 50% (R)"
$ws.Range("I4").Value = "external.R"

# Add new row 5
$ws.Range("A5").Value = "z0bug.wt_1040-23A"
$ws.Range("B5").Value = "1040-23%A"
$ws.Range("C5").Value = "1040 – 23% su 100% (A)"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "z0bug.coa_153110"
$ws.Range("F5").Value = "z0bug.coa_260110"
$ws.Range("G5").Value = "z0bug.jou_misc"
$ws.Range("H5").Value = "account.account_payment_term_15days"
$ws.Range("I5").Value = "external.A"
$ws.Range("J5").Value = "ritenuta"
$ws.Range("K5").Value = 1

# Match the cell style used in column G of the other data rows
$ws.Range("G5").Style = $ws.Range("G2").Style

$ws.Range("A2").Select()
